# ---------------------------------------------------------------------------
# QCA-AID-Explorer-Config.xlsx : add two "Sentiment" analysis config sheets
# (Sentiment1 / Sentiment2), tweak the Heatmap1 "figsize" row, and nudge a
# couple of selections / the active tab -- per the commit "add sentiment
# analysis to Explorer".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Heatmap1: figsize row (row 12) gets a text value "14x10" plus a note
#    "Dimensionen der Grafik" (was a bare number 14.1 with no note before).
#    NB: we write the note (C12) before the value (B12) so that the brand
#    new shared-strings end up in the same relative order as upstream.
# ---------------------------------------------------------------------------
$wsHeatmap = $wb.Worksheets.Item("Heatmap1")
$wsHeatmap.Range("C12").Value = "Dimensionen der Grafik"
$wsHeatmap.Range("B12").Value = "14x10"

# ---------------------------------------------------------------------------
# 2) Build "Sentiment1" as a copy of BenutzerdefinierteAnalyse (same column
#    widths / header styling / row styling), then overwrite its parameter
#    rows for a sentiment analysis of "Akteure" (actors).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("BenutzerdefinierteAnalyse")
$template.Copy($null, $template)
$s1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s1.Name = "Sentiment1"

# row 7: analysis_type -> sentiment_analysis
$s1.Range("B7").Value = "sentiment_analysis"
# row 8: text_column -> "Text", note -> Standard "Text"
$s1.Range("B8").Value = "Text"
$s1.Range("C8").Value = "Auszuwertende Spalte, Standard ""Text"""
# row 9 (new): sentiment_categories
$s1.Range("A9").Value = "sentiment_categories"
$s1.Range("B9").Value = "Positiv, Negativ, Neutral"
# row 10 (new): color_mapping
$s1.Range("A10").Value = "color_mapping"
$s1.Range("B10").Value = '{"Positiv": "#4CAF50", "Negativ": "#F44336", "Neutral": "#9E9E9E"}'
# row 11 (new): chart_title
$s1.Range("A11").Value = "chart_title"
$s1.Range("B11").Value = "Sentiment-Analyse: Akteure (Positiv/Negativ)"
# row 12 (new): temperature
$s1.Range("A12").Value = "temperature"
$s1.Range("B12").Value = "0.3"
# row 13 (new): crosstab_dimensions
$s1.Range("A13").Value = "crosstab_dimensions"
$s1.Range("B13").Value = "Dokument, Hauptkategorie"
# row 14 (new): figsize
$s1.Range("A14").Value = "figsize "
$s1.Range("B14").Value = "12x8"

# --- formatting: stretch the "middle" style (row 8: s3/s3/s4) down over
#     rows 9-13, then restore the "last row" style (row 9: s3/s3/s3) on the
#     very last row (14).
$s1.Range("A9:C9").Copy()
$s1.Range("A14:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$s1.Range("A8:C8").Copy()
$s1.Range("A9:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# re-run the value writes that got clobbered by the format paste above
$s1.Range("A9").Value = "sentiment_categories"
$s1.Range("A10").Value = "color_mapping"
$s1.Range("A11").Value = "chart_title"
$s1.Range("A12").Value = "temperature"
$s1.Range("A13").Value = "crosstab_dimensions"
$s1.Range("A14").Value = "figsize "

# row heights: 9 back to default, 10/11 wrap onto two lines
$s1.Rows.Item(9).RowHeight = 15
$s1.Rows.Item(10).RowHeight = 30
$s1.Rows.Item(11).RowHeight = 30

$s1.Range("F11").Select()

# ---------------------------------------------------------------------------
# 3) Build "Sentiment2" the same way, for a sentiment analysis of
#    "Ressourcen" (resources), including a custom prompt_template.
# ---------------------------------------------------------------------------
$template2 = $wb.Worksheets.Item("BenutzerdefinierteAnalyse")
$template2.Copy($null, $s1)
$s2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s2.Name = "Sentiment2"

$s2.Range("B7").Value = "sentiment_analysis"
$s2.Range("B8").Value = "Text"
$s2.Range("C8").Value = "Auszuwertende Spalte"
$s2.Range("A9").Value = "sentiment_categories"
$s2.Range("B9").Value = "Kritisch, Befürwortend, Ambivalent, Neutral"
$s2.Range("A10").Value = "color_mapping"
$s2.Range("B10").Value = '{"Kritisch": "#FF5722", "Befürwortend": "#2196F3", "Ambivalent": "#9C27B0", "Neutral": "#9E9E9E"}'
$s2.Range("A11").Value = "chart_title"
$s2.Range("B11").Value = "Sentiment-Analyse: Ressourcen (Kritisch/Befürwortend)"
$s2.Range("A12").Value = "temperature"
$s2.Range("B12").Value = "0.3"
$s2.Range("A13").Value = "prompt_template"
$promptRessourcen = @"
Du bist ein Experte für qualitative Textanalyse.
Klassifiziere den folgenden Text in Bezug auf Ressourcen anhand des Sentiments in eine der folgenden Kategorien: Kritisch, Befürwortend, Ambivalent, Neutral
Beachte bei deiner Analyse:
1. Die Bewertung von Ressourcen im Text
2. Die Tonalität gegenüber Ressourcenverfügbarkeit und -nutzung
3. Die impliziten und expliziten Wertungen
Text:
---
[Text kommt hier]
---
Antworte mit einem JSON-Objekt im folgenden Format:
{
    "sentiment": "Kategorie", // Eine der vorgegebenen Kategorien
    "keywords": ["wort1", "wort2", "wort3"], // 3-5 Schlüsselwörter, die zur Bewertung von Ressourcen im Text entscheidend sind
    "explanation": "Kurze Begründung" // Kurze Erklärung (1-2 Sätze)
}
"@
$s2.Range("B13").Value = $promptRessourcen
$s2.Range("C13").Value = "Custom Prompt"
$s2.Range("A14").Value = "crosstab_dimensions"
$s2.Range("B14").Value = "Dokument, Hauptkategorie"
$s2.Range("A15").Value = "figsize "
$s2.Range("B15").Value = "12x8"

$s2.Range("A9:C9").Copy()
$s2.Range("A15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$s2.Range("A8:C8").Copy()
$s2.Range("A9:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$s2.Range("A9").Value = "sentiment_categories"
$s2.Range("A10").Value = "color_mapping"
$s2.Range("A11").Value = "chart_title"
$s2.Range("A12").Value = "temperature"
$s2.Range("A13").Value = "prompt_template"
$s2.Range("A14").Value = "crosstab_dimensions"
$s2.Range("A15").Value = "figsize "

$s2.Rows.Item(9).RowHeight = 30
$s2.Rows.Item(10).RowHeight = 60
$s2.Rows.Item(11).RowHeight = 30
$s2.Rows.Item(12).RowHeight = 15
$s2.Rows.Item(13).RowHeight = 409.5

$s2.Range("I9").Select()

# ---------------------------------------------------------------------------
# 4) View-state touch-ups on the pre-existing sheets.
# ---------------------------------------------------------------------------
$wsHeatmap.Activate()
$wsHeatmap.Range("B13").Select()

$wsCustom = $wb.Worksheets.Item("BenutzerdefinierteAnalyse")
$wsCustom.Activate()
$wsCustom.Range("B9").Select()

# Sentiment2 ends up the active/selected tab, matching the saved workbook.
$s2.Activate()
